{"js": "// Replace the date and the arithmetic problems with new values.\nconst replacements = [\n  [\"2024-06-03 Monday\", \"2024-06-04 Tuesday\"],\n  [\"187\u00d75=\", \"688\u00d78=\"],\n  [\"109\u00d72=\", \"306\u00d75=\"],\n  [\"795\u00d73=\", \"245\u00d72=\"],\n  [\"145\u00d73=\", \"931\u00d73=\"],\n  [\"226\u00d77=\", \"847\u00d75=\"],\n  [\"132\u00d76=\", \"135\u00d73=\"],\n  [\"259\u00d74=\", \"945\u00d72=\"],\n  [\"406\u00d76=\", \"709\u00d73=\"],\n  [\"592\u00d74=\", \"246\u00d79=\"],\n  [\"234\u00d77=\", \"503\u00d74=\"],\n  [\"139\u00d74=\", \"487\u00d77=\"],\n  [\"938\u00d77=\", \"375\u00d75=\"],\n  [\"113\u00d77=\", \"635\u00d76=\"],\n  [\"834\u00d78=\", \"593\u00d76=\"],\n  [\"481\u00d76=\", \"105\u00d78=\"],\n  [\"785\u00d79=\", \"191\u00d79=\"],\n  [\"658\u00d78=\", \"848\u00d79=\"],\n  [\"406\u00d78=\", \"521\u00d79=\"],\n  [\"668\u00d74=\", \"101\u00d73=\"],\n  [\"263\u00d75=\", \"240\u00d77=\"],\n  [\"893\u00d77=\", \"103\u00d77=\"],\n  [\"812\u00d74=\", \"930\u00d72=\"],\n  [\"351\u00d77=\", \"545\u00d74=\"],\n  [\"339\u00d72=\", \"930\u00d73=\"],\n  [\"785\u00d72=\", \"563\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-06-03 Monday\"; New = \"2024-06-04 Tuesday\" },\n    @{ Old = \"187\u00d75=\"; New = \"688\u00d78=\" },\n    @{ Old = \"109\u00d72=\"; New = \"306\u00d75=\" },\n    @{ Old = \"795\u00d73=\"; New = \"245\u00d72=\" },\n    @{ Old = \"145\u00d73=\"; New = \"931\u00d73=\" },\n    @{ Old = \"226\u00d77=\"; New = \"847\u00d75=\" },\n    @{ Old = \"132\u00d76=\"; New = \"135\u00d73=\" },\n    @{ Old = \"259\u00d74=\"; New = \"945\u00d72=\" },\n    @{ Old = \"406\u00d76=\"; New = \"709\u00d73=\" },\n    @{ Old = \"592\u00d74=\"; New = \"246\u00d79=\" },\n    @{ Old = \"234\u00d77=\"; New = \"503\u00d74=\" },\n    @{ Old = \"139\u00d74=\"; New = \"487\u00d77=\" },\n    @{ Old = \"938\u00d77=\"; New = \"375\u00d75=\" },\n    @{ Old = \"113\u00d77=\"; New = \"635\u00d76=\" },\n    @{ Old = \"834\u00d78=\"; New = \"593\u00d76=\" },\n    @{ Old = \"481\u00d76=\"; New = \"105\u00d78=\" },\n    @{ Old = \"785\u00d79=\"; New = \"191\u00d79=\" },\n    @{ Old = \"658\u00d78=\"; New = \"848\u00d79=\" },\n    @{ Old = \"406\u00d78=\"; New = \"521\u00d79=\" },\n    @{ Old = \"668\u00d74=\"; New = \"101\u00d73=\" },\n    @{ Old = \"263\u00d75=\"; New = \"240\u00d77=\" },\n    @{ Old = \"893\u00d77=\"; New = \"103\u00d77=\" },\n    @{ Old = \"812\u00d74=\"; New = \"930\u00d72=\" },\n    @{ Old = \"351\u00d77=\"; New = \"545\u00d74=\" },\n    @{ Old = \"339\u00d72=\"; New = \"930\u00d73=\" },\n    @{ Old = \"785\u00d72=\"; New = \"563\u00d76=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
